$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add header cells I1 / J1, copying the formatting of the existing H1 header ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Fill in the I (I0) and J (IF) data columns for rows 2-64 ---
$iVals = @(7,8,7,7,8,9,5,5,7,6,9,7,8,9,9,11,9,9,9,9,8,8,9,8,9,9,9,9,7,8,9,10,9,9,8,8,9,9,9,8,9,9,9,9,9,9,8,9,9,9,9,9,9,7,8,8,8,7,3,7,8,6,6)
$jVals = @(7,8,7,7,8,9,6,5,7,6,9,8,8,9,9,11,9,9,9,9,8,8,9,9,9,9,9,9,8,8,9,10,9,9,8,8,10,9,9,9,9,9,9,10,9,9,9,9,9,9,10,9,9,8,8,8,8,7,4,8,8,6,6)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
  $row = $idx + 2
  $ws.Cells.Item($row, 9).Value = $iVals[$idx]
  $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}

Write-Host "I0/IF columns populated"
